$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet dimension implicitly grows as we write data; Excel COM
# recalculates the used range automatically once new cells are populated.

$data = @(
    ,@("Menu", 17, 5, 22)
    ,@("DataGrid", 5, 6, 11)
    ,@("Combobox", 6, 4, 10)
    ,@("Nav", 8, 1, 9)
    ,@("Tree", 9, 0, 9)
    ,@("Popover", 6, 2, 8)
    ,@("Dialog", 5, 1, 6)
    ,@("Dropdown", 4, 2, 6)
    ,@("Table", 4, 1, 5)
    ,@("Tooltip", 4, 1, 5)
    ,@("Virtualizer", 4, 0, 4)
    ,@("TagPicker", 3, 1, 4)
    ,@("Toolbar", 3, 1, 4)
    ,@("Tabs", 2, 1, 3)
    ,@("Skeleton", 2, 1, 3)
    ,@("TeachingPopover", 3, 0, 3)
    ,@("Calendar Compat", 3, 0, 3)
    ,@("Drawer", 2, 0, 2)
    ,@("MessageBar", 2, 0, 2)
    ,@("FluentProvider", 0, 2, 2)
    ,@("Slider", 0, 2, 2)
    ,@("Toast", 1, 1, 2)
    ,@("List", 2, 0, 2)
    ,@("Accordion", 2, 0, 2)
    ,@("DatePicker", 0, 2, 2)
    ,@("Switch", 2, 0, 2)
    ,@("Portal", 2, 0, 2)
    ,@("Carousel", 1, 0, 1)
    ,@("AvatarGroup", 1, 0, 1)
    ,@("SplitButton", 1, 0, 1)
    ,@("SearchBox", 1, 0, 1)
    ,@("Label", 1, 0, 1)
    ,@("InfoLabel", 1, 0, 1)
    ,@("Spinner", 0, 1, 1)
    ,@("Button", 0, 1, 1)
    ,@("Persona", 1, 0, 1)
    ,@("FocusTrapZone", 1, 0, 1)
    ,@("Input", 0, 1, 1)
    ,@("Checkbox", 1, 0, 1)
    ,@("Avatar", 1, 0, 1)
    ,@("Popup", 1, 0, 1)
    ,@("DatePickerCompat", 0, 1, 1)
    ,@("Image", 0, 1, 1)
    ,@("MenuItem", 1, 0, 1)
    ,@("Badge", 1, 0, 1)
    ,@("Tag", 0, 0, 0)
    ,@("SpinButton", 0, 0, 0)
    ,@("Segment", 0, 0, 0)
    ,@("Rating", 0, 0, 0)
    ,@("Keytip", 0, 0, 0)
    ,@("Pickers", 0, 0, 0)
    ,@("ColorPicker", 0, 0, 0)
    ,@("FloatingLabelInput", 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

